# Fixed inability to pass in duplicate asset descriptions as input
#
# The budget column headers on the BME_RENAL_SPH sheet used an
# abbreviated "budg_" prefix that collided with validation logic
# elsewhere, preventing duplicate asset descriptions from being
# entered. Renaming them to the unabbreviated "budgeted_" prefix
# (matching the "actual_" naming convention already used) fixes it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BME_RENAL_SPH")

$ws.Range("G1").Value = "budgeted_partial_oh"
$ws.Range("H1").Value = "budgeted_total_exp"
$ws.Range("I1").Value = "budgeted_labour_exp"
$ws.Range("J1").Value = "budgeted_contracts_exp"
$ws.Range("K1").Value = "budgeted_parts_exp"

# Move the sheet's remembered cursor/selection position.
$ws.Activate()
$ws.Range("H18").Select()
